# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" sheets to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 268
    3  = 284
    4  = 292
    5  = 845
    7  = 303
    8  = 8259
    9  = 76
    12 = 110
    15 = 22
    18 = 252
    19 = 712
    20 = 29
    21 = 81
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
